# Updated symbol list on Fri Jan 13 22:25:32 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures on the
# "cryptos" sheet with the latest snapshot pulled by the scraper.
#
# The sheet stores Price/Volume as literal text (e.g. "292.21", "2.04%")
# rather than numbers, so every new value is written with a leading
# apostrophe. That keeps Excel from auto-coercing a numeric-looking
# string ("292.33") into a Number, or a "x.xx%" string into a Percentage,
# which would silently change the cell's stored type/precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = 4; Text = "292.33" },
    @{ Row = 2; Col = 5; Text = "2.09%" },
    @{ Row = 3; Col = 4; Text = "29.67" },
    @{ Row = 3; Col = 5; Text = "3.83%" },
    @{ Row = 4; Col = 4; Text = "5.282" },
    @{ Row = 4; Col = 5; Text = "4.23%" },
    @{ Row = 5; Col = 4; Text = "0.07173" },
    @{ Row = 5; Col = 5; Text = "8.03%" },
    @{ Row = 6; Col = 4; Text = "7.537" },
    @{ Row = 6; Col = 5; Text = "2.31%" },
    @{ Row = 7; Col = 4; Text = "3.594" },
    @{ Row = 7; Col = 5; Text = "5.91%" },
    @{ Row = 8; Col = 4; Text = "1.399" },
    @{ Row = 8; Col = 5; Text = "2.21%" },
    @{ Row = 9; Col = 4; Text = "0.9085" },
    @{ Row = 9; Col = 5; Text = "-3.18%" },
    @{ Row = 10; Col = 4; Text = "0.1621" },
    @{ Row = 10; Col = 5; Text = "3.57%" },
    @{ Row = 11; Col = 4; Text = "0.07626" },
    @{ Row = 11; Col = 5; Text = "15.57%" },
    @{ Row = 12; Col = 4; Text = "0.07782" },
    @{ Row = 12; Col = 5; Text = "2.38%" },
    @{ Row = 13; Col = 4; Text = "0.02913" },
    @{ Row = 13; Col = 5; Text = "-0.79%" },
    @{ Row = 14; Col = 4; Text = "0.09000" },
    @{ Row = 14; Col = 5; Text = "0.10%" },
    @{ Row = 15; Col = 4; Text = "0.001583" },
    @{ Row = 15; Col = 5; Text = "0.07%" },
    @{ Row = 16; Col = 4; Text = "0.0006529" },
    @{ Row = 16; Col = 5; Text = "0.78%" },
    @{ Row = 17; Col = 4; Text = "0.006114" },
    @{ Row = 17; Col = 5; Text = "-2.52%" },
    @{ Row = 18; Col = 4; Text = "3.478" },
    @{ Row = 18; Col = 5; Text = "1.10%" },
    @{ Row = 19; Col = 4; Text = "2.233" },
    @{ Row = 19; Col = 5; Text = "-0.94%" },
    @{ Row = 20; Col = 4; Text = "0.3252" },
    @{ Row = 20; Col = 5; Text = "1.11%" },
    @{ Row = 21; Col = 4; Text = "0.1366" },
    @{ Row = 21; Col = 5; Text = "5.22%" },
    @{ Row = 22; Col = 4; Text = "4.038" },
    @{ Row = 22; Col = 5; Text = "-1.04%" },
    @{ Row = 23; Col = 4; Text = "0.1591" },
    @{ Row = 23; Col = 5; Text = "2.43%" },
    @{ Row = 24; Col = 5; Text = "0.56%" },
    @{ Row = 25; Col = 5; Text = "1.98%" },
    @{ Row = 26; Col = 4; Text = "0.004249" },
    @{ Row = 26; Col = 5; Text = "-4.91%" },
    @{ Row = 27; Col = 4; Text = "0.0001164" },
    @{ Row = 27; Col = 5; Text = "-6.92%" },
    @{ Row = 28; Col = 4; Text = "0.0001681" },
    @{ Row = 28; Col = 5; Text = "3.84%" },
    @{ Row = 40; Col = 4; Text = "0.04432" },
    @{ Row = 40; Col = 5; Text = "5.41%" },
    @{ Row = 41; Col = 4; Text = "0.007000" },
    @{ Row = 41; Col = 5; Text = "3.73%" },
    @{ Row = 42; Col = 4; Text = "0.1276" },
    @{ Row = 42; Col = 5; Text = "2.21%" },
    @{ Row = 43; Col = 4; Text = "0.002198" },
    @{ Row = 43; Col = 5; Text = "8.80%" },
    @{ Row = 44; Col = 4; Text = "0.01320" },
    @{ Row = 44; Col = 5; Text = "7.80%" },
    @{ Row = 45; Col = 4; Text = "0.00005837" },
    @{ Row = 45; Col = 5; Text = "4.18%" },
    @{ Row = 47; Col = 4; Text = "0.01293" },
    @{ Row = 47; Col = 5; Text = "-1.11%" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = "'" + $u.Text
}
